$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear all existing content so we can rebuild the table from scratch with
# the new "_requirements" column inserted between "id" and "withholding_tax_id".
$ws.Cells.Clear()

# --- Header row -----------------------------------------------------------
$ws.Cells.Item(1,1).Value = "id"
$ws.Cells.Item(1,2).Value = "_requirements"
$ws.Cells.Item(1,3).Value = "withholding_tax_id"
$ws.Cells.Item(1,4).Value = "tax"
$ws.Cells.Item(1,5).Value = "base"

# --- Row 2: z0bug.wt_1040_1 ------------------------------------------------
$ws.Cells.Item(2,1).Value = "z0bug.wt_1040_1"
$ws.Cells.Item(2,3).Value = "z0bug.wt_1040"
$ws.Cells.Item(2,4).Value = 20
$ws.Cells.Item(2,5).Value = 1

# --- Row 3: z0bug.wt_1038_1 (zero case) ------------------------------------
$ws.Cells.Item(3,1).Value = "z0bug.wt_1038_1"
$ws.Cells.Item(3,2).Value = "G=='zero'"
$ws.Cells.Item(3,3).Value = "z0bug.wt_1038"
$ws.Cells.Item(3,4).Value = 23
$ws.Cells.Item(3,5).NumberFormat = "@"
$ws.Cells.Item(3,5).Value = "0.5"
$ws.Cells.Item(3,5).NumberFormat = "General"

# --- Row 4: z0bug.wt_1038_1 (non-zero case) --------------------------------
$ws.Cells.Item(4,1).Value = "z0bug.wt_1038_1"
$ws.Cells.Item(4,2).Value = "G!='zero'"
$ws.Cells.Item(4,3).Value = "z0bug.wt_1038"
$ws.Cells.Item(4,4).NumberFormat = "@"
$ws.Cells.Item(4,4).Value = "11.5"
$ws.Cells.Item(4,4).NumberFormat = "General"
$ws.Cells.Item(4,5).Value = 1

# --- Row 5: z0bug.wt_enasarco_1_1 (zero case) ------------------------------
$ws.Cells.Item(5,1).Value = "z0bug.wt_enasarco_1_1"
$ws.Cells.Item(5,2).Value = "G=='zero'"
$ws.Cells.Item(5,3).Value = "z0bug.wt_enasarco_1"
$ws.Cells.Item(5,4).Value = 17
$ws.Cells.Item(5,5).NumberFormat = "@"
$ws.Cells.Item(5,5).Value = "0.5"
$ws.Cells.Item(5,5).NumberFormat = "General"

# --- Row 6: z0bug.wt_enasarco_1_1 (non-zero case) --------------------------
$ws.Cells.Item(6,1).Value = "z0bug.wt_enasarco_1_1"
$ws.Cells.Item(6,2).Value = "G!='zero'"
$ws.Cells.Item(6,3).Value = "z0bug.wt_enasarco_1"
$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value = "8.5"
$ws.Cells.Item(6,4).NumberFormat = "General"
$ws.Cells.Item(6,5).Value = 1

# --- Row 7: z0bug.wt_1040-23A_1 --------------------------------------------
$ws.Cells.Item(7,1).Value = "z0bug.wt_1040-23A_1"
$ws.Cells.Item(7,3).Value = "z0bug.wt_1040-23A"
$ws.Cells.Item(7,4).Value = 23
$ws.Cells.Item(7,5).Value = 1

# --- Column widths ----------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 14.0
$ws.Columns.Item(3).ColumnWidth = 18.667
$ws.Columns.Item(4).ColumnWidth = 3.167
$ws.Columns.Item(5).ColumnWidth = 4.667

# --- Selection --------------------------------------------------------------
$ws.Range("D5").Select()
